# Auto-generated edit script
$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (exhibitions) - F column "want to go" count updates ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value = 393
$ws1.Cells.Item(6, 6).Value = 1254
$ws1.Cells.Item(8, 6).Value = 101
$ws1.Cells.Item(9, 6).Value = 206
$ws1.Cells.Item(11, 6).Value = 185
$ws1.Cells.Item(15, 6).Value = 200
$ws1.Cells.Item(16, 6).Value = 1531
$ws1.Cells.Item(19, 6).Value = 358
$ws1.Cells.Item(21, 6).Value = 850
$ws1.Cells.Item(26, 6).Value = 1470
$ws1.Cells.Item(27, 6).Value = 70
$ws1.Cells.Item(28, 6).Value = 50
$ws1.Cells.Item(29, 6).Value = 449
$ws1.Cells.Item(30, 6).Value = 715
$ws1.Cells.Item(31, 6).Value = 1328
$ws1.Cells.Item(33, 6).Value = 1418
$ws1.Cells.Item(38, 6).Value = 692
$ws1.Cells.Item(39, 6).Value = 874
$ws1.Cells.Item(41, 6).Value = 261

# ---- Sheet: 演出 (performances) - F column "want to go" count updates ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(6, 6).Value = 181
$ws2.Cells.Item(10, 6).Value = 4
$ws2.Cells.Item(15, 6).Value = 651

# ---- Sheet: 全部类型 (all types) - rows 3-35 replaced with updated/re-sorted merged list ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 2).Value = "2024-05-12"
$ws4.Cells.Item(3, 3).Value = "杭州·音乐番ONLY"
$ws4.Cells.Item(3, 4).Value = "体育场路武林广场11号 杭州大厦中央商城"
$ws4.Cells.Item(3, 5).Value = "2024.05.12 10:00-05.12 16:00"
$ws4.Cells.Item(3, 6).Value = 140
$ws4.Cells.Item(3, 7).Value = 58
$ws4.Cells.Item(3, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84656"
$ws4.Cells.Item(3, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/eauUzJj11713741020828.jpeg"

$ws4.Cells.Item(4, 2).Value = "2024-05-18"
$ws4.Cells.Item(4, 3).Value = "杭州·HCCL·高校联展"
$ws4.Cells.Item(4, 4).Value = "康候圣街99号 顺丰创新中心"
$ws4.Cells.Item(4, 5).Value = "2024.05.18 09:00-05.19 19:00"
$ws4.Cells.Item(4, 6).Value = 393
$ws4.Cells.Item(4, 7).Value = 49
$ws4.Cells.Item(4, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84561"
$ws4.Cells.Item(4, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/uyBT10rf1713186863701.png"

$ws4.Cells.Item(5, 2).Value = "2024-05-18"
$ws4.Cells.Item(5, 3).Value = "杭州·Jo迪"
$ws4.Cells.Item(5, 4).Value = "萧杭路28号 格拉斯club"
$ws4.Cells.Item(5, 5).Value = "2024.05.18 13:00-05.18 19:00"
$ws4.Cells.Item(5, 6).Value = 200
$ws4.Cells.Item(5, 7).Value = 198
$ws4.Cells.Item(5, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83008"
$ws4.Cells.Item(5, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/AEtl5BHN1711015003341.jpeg"

$ws4.Cells.Item(6, 2).Value = "2024-05-18"
$ws4.Cells.Item(6, 3).Value = "杭州·SK11三坑特卖会"
$ws4.Cells.Item(6, 4).Value = "鸿泰路与明月桥路交汇处东南角方位(杭港地铁1号线/杭州地铁4号线彭埠站D口20米) 港龙悠乐城"
$ws4.Cells.Item(6, 5).Value = "2024.05.18 10:00-05.19 19:00"
$ws4.Cells.Item(6, 6).Value = 33
$ws4.Cells.Item(6, 7).Value = 99
$ws4.Cells.Item(6, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85148"
$ws4.Cells.Item(6, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/lio1sGir1714982499357.jpeg"

$ws4.Cells.Item(7, 2).Value = "2024-05-18"
$ws4.Cells.Item(7, 3).Value = "杭州·《天空之城》久石让·宫崎骏动漫经典作品音乐会·筑乐之城"
$ws4.Cells.Item(7, 4).Value = "曙光路31号 浙江音乐厅"
$ws4.Cells.Item(7, 5).Value = "2024.05.18 19:30-05.18 21:00"
$ws4.Cells.Item(7, 6).Value = 7
$ws4.Cells.Item(7, 7).Value = 100
$ws4.Cells.Item(7, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84905"
$ws4.Cells.Item(7, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/OocULytV1714103273912.jpeg"

$ws4.Cells.Item(8, 2).Value = "2024-05-18"
$ws4.Cells.Item(8, 3).Value = "杭州·《沐云华·次元狂想》经典动漫二次元ACG音乐会"
$ws4.Cells.Item(8, 4).Value = "建国南路280号 杭州红星剧院"
$ws4.Cells.Item(8, 5).Value = "2024.05.18 19:30-05.18 22:00"
$ws4.Cells.Item(8, 6).Value = 181
$ws4.Cells.Item(8, 7).Value = 252
$ws4.Cells.Item(8, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83113"
$ws4.Cells.Item(8, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/TXmgAvCC1710582339525.jpeg"

$ws4.Cells.Item(9, 2).Value = "2024-05-18"
$ws4.Cells.Item(9, 3).Value = "杭州·现世繁华-代号鸢only"
$ws4.Cells.Item(9, 4).Value = "丁城路丁桥桃花湖公园北区 典酷沉浸式艺术空间(桃花湖店)"
$ws4.Cells.Item(9, 5).Value = "2024.05.18 10:00-05.18 21:00"
$ws4.Cells.Item(9, 6).Value = 1254
$ws4.Cells.Item(9, 7).Value = 156
$ws4.Cells.Item(9, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81905"
$ws4.Cells.Item(9, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/m3upuV2F1708327958926.jpeg"

$ws4.Cells.Item(10, 2).Value = "2024-05-18"
$ws4.Cells.Item(10, 3).Value = "杭州·第五幼儿园·第五人格only展"
$ws4.Cells.Item(10, 4).Value = "康候圣街99号 顺丰创新中心"
$ws4.Cells.Item(10, 5).Value = "2024.05.18 09:00-05.18 17:00"
$ws4.Cells.Item(10, 6).Value = 459
$ws4.Cells.Item(10, 7).Value = 68
$ws4.Cells.Item(10, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82834"
$ws4.Cells.Item(10, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/uum0yj2L1713577581499.jpeg"

$ws4.Cells.Item(11, 2).Value = "2024-05-19"
$ws4.Cells.Item(11, 3).Value = "杭州·m字刘海少年和粉毛少女only"
$ws4.Cells.Item(11, 4).Value = "康候圣街99号 顺丰创新中心"
$ws4.Cells.Item(11, 5).Value = "2024.05.19 09:00-05.19 17:00"
$ws4.Cells.Item(11, 6).Value = 101
$ws4.Cells.Item(11, 7).Value = 68
$ws4.Cells.Item(11, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82831"
$ws4.Cells.Item(11, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/bVvk6Eky1710383662942.jpeg"

$ws4.Cells.Item(12, 2).Value = "2024-05-19"
$ws4.Cells.Item(12, 3).Value = "杭州·原X铁X崩ONLY"
$ws4.Cells.Item(12, 4).Value = "黄姑山路51-4号 0101park"
$ws4.Cells.Item(12, 5).Value = "2024.05.19 10:00-05.19 17:30"
$ws4.Cells.Item(12, 6).Value = 206
$ws4.Cells.Item(12, 7).Value = 58
$ws4.Cells.Item(12, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84849"
$ws4.Cells.Item(12, 9).Value = "//i2.hdslb.com/bfs/openplatform/202405/pfNAP8zt1715074013459.jpeg"

$ws4.Cells.Item(13, 2).Value = "2024-05-19"
$ws4.Cells.Item(13, 3).Value = "杭州·第二届动漫迷城嘉年华"
$ws4.Cells.Item(13, 4).Value = "体育场路武林广场11号 杭州大厦中央商城"
$ws4.Cells.Item(13, 5).Value = "2024.05.19 10:00-05.19 17:00"
$ws4.Cells.Item(13, 6).Value = 159
$ws4.Cells.Item(13, 7).Value = 60
$ws4.Cells.Item(13, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83964"
$ws4.Cells.Item(13, 9).Value = "//i1.hdslb.com/bfs/openplatform/202404/3WNfBWY61713863269103.jpeg"

$ws4.Cells.Item(14, 2).Value = "2024-05-25"
$ws4.Cells.Item(14, 3).Value = "杭州·D3动漫游戏嘉年华"
$ws4.Cells.Item(14, 4).Value = "德胜东路2539号 梦马汽车小镇"
$ws4.Cells.Item(14, 5).Value = "2024.05.25 10:00-05.25 17:00"
$ws4.Cells.Item(14, 6).Value = 185
$ws4.Cells.Item(14, 7).Value = 50
$ws4.Cells.Item(14, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84912"
$ws4.Cells.Item(14, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/HXY7pTYI1715053764601.jpeg"

$ws4.Cells.Item(15, 2).Value = "2024-05-25"
$ws4.Cells.Item(15, 3).Value = "杭州·Redamancy动漫游戏嘉年华×运动番全明星"
$ws4.Cells.Item(15, 4).Value = "富春路80号(甬江路地铁站A口旁) 杭州全民健身中心"
$ws4.Cells.Item(15, 5).Value = "2024.05.25 10:00-05.26 17:00"
$ws4.Cells.Item(15, 6).Value = 1060
$ws4.Cells.Item(15, 7).Value = 68
$ws4.Cells.Item(15, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84947"
$ws4.Cells.Item(15, 9).Value = "//i1.hdslb.com/bfs/openplatform/202404/65Usx6BT1713796309433.jpeg"

$ws4.Cells.Item(16, 2).Value = "2024-05-25"
$ws4.Cells.Item(16, 3).Value = "杭州·原神X星铁X绝区零only"
$ws4.Cells.Item(16, 4).Value = "望江东路333号 杭州瑞莱克斯大酒店"
$ws4.Cells.Item(16, 5).Value = "2024.05.25 10:00-05.25 17:00"
$ws4.Cells.Item(16, 6).Value = 273
$ws4.Cells.Item(16, 7).Value = 60
$ws4.Cells.Item(16, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82754"
$ws4.Cells.Item(16, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/qA0LNJuF1710234461030.jpeg"

$ws4.Cells.Item(17, 2).Value = "2024-05-25"
$ws4.Cells.Item(17, 3).Value = "杭州·早鸟5折起·《LALALAND爱乐之城》浪漫主题音乐会"
$ws4.Cells.Item(17, 4).Value = "武林路77号 浙江省文化馆小剧场（原群艺馆小剧场）"
$ws4.Cells.Item(17, 5).Value = "2024.05.25 19:30-05.25 21:00"
$ws4.Cells.Item(17, 6).Value = 9
$ws4.Cells.Item(17, 7).Value = 100
$ws4.Cells.Item(17, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84519"
$ws4.Cells.Item(17, 9).Value = "//i1.hdslb.com/bfs/openplatform/202404/jJLft5tT1712888683239.jpeg"

$ws4.Cells.Item(18, 2).Value = "2024-05-25"
$ws4.Cells.Item(18, 3).Value = "杭州·热血番ONLY1.0"
$ws4.Cells.Item(18, 4).Value = "康候圣街99号 顺丰创新中心"
$ws4.Cells.Item(18, 5).Value = "2024.05.25 10:00-05.26 17:00"
$ws4.Cells.Item(18, 6).Value = 200
$ws4.Cells.Item(18, 7).Value = 68
$ws4.Cells.Item(18, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85042"
$ws4.Cells.Item(18, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/HSZzsbLs1714221391790.jpeg"

$ws4.Cells.Item(19, 2).Value = "2024-05-25"
$ws4.Cells.Item(19, 3).Value = "杭州·第三届缘起cp展 我们二次元的情人节！"
$ws4.Cells.Item(19, 4).Value = "黄姑山路51-4号 0101park"
$ws4.Cells.Item(19, 5).Value = "2024.05.25 10:00-05.26 17:00"
$ws4.Cells.Item(19, 6).Value = 1531
$ws4.Cells.Item(19, 7).Value = 65
$ws4.Cells.Item(19, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83336"
$ws4.Cells.Item(19, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/D9t8ms7G1711350634757.png"

$ws4.Cells.Item(20, 2).Value = "2024-05-26"
$ws4.Cells.Item(20, 3).Value = "杭州·恋与深空×恋与制作人only"
$ws4.Cells.Item(20, 4).Value = "望江东路333号 杭州瑞莱克斯大酒店"
$ws4.Cells.Item(20, 5).Value = "2024.05.26 10:00-05.26 17:00"
$ws4.Cells.Item(20, 6).Value = 560
$ws4.Cells.Item(20, 7).Value = 60
$ws4.Cells.Item(20, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84077"
$ws4.Cells.Item(20, 9).Value = "//i1.hdslb.com/bfs/openplatform/202404/V6V4Pppv1712736555042.jpeg"

$ws4.Cells.Item(21, 2).Value = "2024-05-26"
$ws4.Cells.Item(21, 3).Value = "杭州·运动番ONLY"
$ws4.Cells.Item(21, 4).Value = "体育场路武林广场11号 杭州大厦中央商城"
$ws4.Cells.Item(21, 5).Value = "2024.05.26 10:00-05.26 16:00"
$ws4.Cells.Item(21, 6).Value = 238
$ws4.Cells.Item(21, 7).Value = 60
$ws4.Cells.Item(21, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84618"
$ws4.Cells.Item(21, 9).Value = "//i1.hdslb.com/bfs/openplatform/202404/anf0T3BA1713742714789.jpeg"

$ws4.Cells.Item(22, 2).Value = "2024-06-01"
$ws4.Cells.Item(22, 3).Value = "杭州·蔚蓝档案only"
$ws4.Cells.Item(22, 4).Value = "北干街道萧杭路689号 杭州时尚外滩艺术中心"
$ws4.Cells.Item(22, 5).Value = "2024.06.01 09:00-06.01 18:00"
$ws4.Cells.Item(22, 6).Value = 358
$ws4.Cells.Item(22, 7).Value = 80
$ws4.Cells.Item(22, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84478"
$ws4.Cells.Item(22, 9).Value = "//i1.hdslb.com/bfs/openplatform/202404/z5lgl4tb1712719299126.jpeg"

$ws4.Cells.Item(23, 2).Value = "2024-06-05"
$ws4.Cells.Item(23, 3).Value = "杭州·英雄时代2024哈瓦西钢琴演奏会"
$ws4.Cells.Item(23, 4).Value = "中国杭州北山路86号西湖岳湖景区 中国杭州西湖岳湖景区印象西湖"
$ws4.Cells.Item(23, 5).Value = "2024.06.05 20:00-06.05 21:30"
$ws4.Cells.Item(23, 6).Value = 2
$ws4.Cells.Item(23, 7).Value = 499
$ws4.Cells.Item(23, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83902"
$ws4.Cells.Item(23, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/BFRFmKpT1712569969076.jpeg"

$ws4.Cells.Item(24, 2).Value = "2024-06-07"
$ws4.Cells.Item(24, 3).Value = "杭州·【鼓楼西戏剧】周一围领衔主演·《枕头人》10周年纪念版"
$ws4.Cells.Item(24, 4).Value = "杭州市江干区新业路39号 杭州大剧院"
$ws4.Cells.Item(24, 5).Value = "2024.06.07 19:30-06.08 22:00"
$ws4.Cells.Item(24, 6).Value = 1
$ws4.Cells.Item(24, 7).Value = 480
$ws4.Cells.Item(24, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84902"
$ws4.Cells.Item(24, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/VZQS8SJP1714020772683.jpeg"

$ws4.Cells.Item(25, 2).Value = "2024-06-08"
$ws4.Cells.Item(25, 3).Value = "杭州·第38届漫展x原崩铁only"
$ws4.Cells.Item(25, 4).Value = "康候圣街99号 顺丰创新中心"
$ws4.Cells.Item(25, 5).Value = "2024.06.08 10:30-06.09 17:00"
$ws4.Cells.Item(25, 6).Value = 1164
$ws4.Cells.Item(25, 7).Value = 60
$ws4.Cells.Item(25, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84802"
$ws4.Cells.Item(25, 9).Value = "//i1.hdslb.com/bfs/openplatform/202404/QsVzW6XP1712908414935.jpeg"

$ws4.Cells.Item(26, 2).Value = "2024-06-09"
$ws4.Cells.Item(26, 3).Value = "杭州·第三届日夜国乙only"
$ws4.Cells.Item(26, 4).Value = "创意路1号 中国智谷富春园区"
$ws4.Cells.Item(26, 5).Value = "2024.06.09 10:00-06.09 23:00"
$ws4.Cells.Item(26, 6).Value = 2682
$ws4.Cells.Item(26, 7).Value = 58
$ws4.Cells.Item(26, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82618"
$ws4.Cells.Item(26, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/fXRzYEFH1710124366279.png"

$ws4.Cells.Item(27, 2).Value = "2024-06-14"
$ws4.Cells.Item(27, 3).Value = "杭州·苗阜王声 青曲社相声全国巡演"
$ws4.Cells.Item(27, 4).Value = "湖墅南路138号 杭州浙话艺术剧院"
$ws4.Cells.Item(27, 5).Value = "2024.06.14 19:30-06.14 22:00"
$ws4.Cells.Item(27, 6).Value = 8
$ws4.Cells.Item(27, 7).Value = 280
$ws4.Cells.Item(27, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83382"
$ws4.Cells.Item(27, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/hUGL3xz01711346789039.jpeg"

$ws4.Cells.Item(28, 2).Value = "2024-06-15"
$ws4.Cells.Item(28, 3).Value = "杭州·次元盛典1.0"
$ws4.Cells.Item(28, 4).Value = "康候圣街99号 顺丰创新中心"
$ws4.Cells.Item(28, 5).Value = "2024.06.15 10:00-06.16 17:00"
$ws4.Cells.Item(28, 6).Value = 1470
$ws4.Cells.Item(28, 7).Value = 68
$ws4.Cells.Item(28, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83672"
$ws4.Cells.Item(28, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/yZAi07mM1712033477653.jpeg"

$ws4.Cells.Item(29, 2).Value = "2024-06-15"
$ws4.Cells.Item(29, 3).Value = "杭州·第三届动漫迷城嘉年华·毕业泳池"
$ws4.Cells.Item(29, 4).Value = "东新路21号 九龙仓君玺"
$ws4.Cells.Item(29, 5).Value = "2024.06.15 10:00-06.15 17:00"
$ws4.Cells.Item(29, 6).Value = 70
$ws4.Cells.Item(29, 7).Value = 70
$ws4.Cells.Item(29, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84338"
$ws4.Cells.Item(29, 9).Value = "//i1.hdslb.com/bfs/openplatform/202404/wQAlXTnK1713202337669.jpeg"

$ws4.Cells.Item(30, 2).Value = "2024-06-22"
$ws4.Cells.Item(30, 3).Value = "杭州·《1+1》2  X PianoLab  “琴梦山川”倪海瑽钢琴独奏音乐会"
$ws4.Cells.Item(30, 4).Value = "杭州市江干区新业路39号 杭州大剧院(音乐厅)"
$ws4.Cells.Item(30, 5).Value = "2024.06.22 19:30-06.22 21:00"
$ws4.Cells.Item(30, 6).Value = 1
$ws4.Cells.Item(30, 7).Value = 100
$ws4.Cells.Item(30, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84521"
$ws4.Cells.Item(30, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/LUu4JB9O1711766011505.jpeg"

$ws4.Cells.Item(31, 2).Value = "2024-06-22"
$ws4.Cells.Item(31, 3).Value = "杭州·巅峰对决·排球少年ONLY"
$ws4.Cells.Item(31, 4).Value = "金桥北路990号 万达广场(杭州富阳店)"
$ws4.Cells.Item(31, 5).Value = "2024.06.22 10:00-06.22 17:00"
$ws4.Cells.Item(31, 6).Value = 50
$ws4.Cells.Item(31, 7).Value = 60
$ws4.Cells.Item(31, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85095"
$ws4.Cells.Item(31, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/3WtpRjjo1714197500930.jpeg"

$ws4.Cells.Item(32, 2).Value = "2024-06-23"
$ws4.Cells.Item(32, 3).Value = "杭州·《亚米·跨越二次元》ACG经典动漫视听音乐会"
$ws4.Cells.Item(32, 4).Value = "金沙大道681号 金沙湖大剧院"
$ws4.Cells.Item(32, 5).Value = "2024.06.23 19:30-06.23 21:10"
$ws4.Cells.Item(32, 6).Value = 23
$ws4.Cells.Item(32, 7).Value = 80
$ws4.Cells.Item(32, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84041"
$ws4.Cells.Item(32, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/UhUuHfad1712564787267.jpeg"

$ws4.Cells.Item(33, 2).Value = "2024-06-23"
$ws4.Cells.Item(33, 3).Value = "杭州·【早鸟5折】中西合奏·再现经典《青城山下·千年等一回》传世国风跨界音乐会"
$ws4.Cells.Item(33, 4).Value = "曙光路31号 浙江音乐厅"
$ws4.Cells.Item(33, 5).Value = "2024.06.23 15:00-06.23 21:00"
$ws4.Cells.Item(33, 6).Value = 4
$ws4.Cells.Item(33, 7).Value = 50
$ws4.Cells.Item(33, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84597"
$ws4.Cells.Item(33, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/jNu5hjYv1713514034369.jpeg"

$ws4.Cells.Item(34, 2).Value = "2024-06-29"
$ws4.Cells.Item(34, 3).Value = "杭州·HD·01"
$ws4.Cells.Item(34, 4).Value = "钱江世纪城奔竞大道353号 杭州国际博览中心"
$ws4.Cells.Item(34, 5).Value = "2024.06.29 10:00-06.30 17:00"
$ws4.Cells.Item(34, 6).Value = 449
$ws4.Cells.Item(34, 7).Value = 75
$ws4.Cells.Item(34, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85006"
$ws4.Cells.Item(34, 9).Value = "//i1.hdslb.com/bfs/openplatform/202404/3ia8Pqym1714378136900.jpeg"

$ws4.Cells.Item(35, 2).Value = "2024-06-29"
$ws4.Cells.Item(35, 3).Value = "杭州·乌托邦次元聚会3.0·二次元全女性夜场"
$ws4.Cells.Item(35, 4).Value = "保淑路2号 The Queen皇后"
$ws4.Cells.Item(35, 5).Value = "2024.06.29 13:00-06.29 19:00"
$ws4.Cells.Item(35, 6).Value = 715
$ws4.Cells.Item(35, 7).Value = 188
$ws4.Cells.Item(35, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84558"
$ws4.Cells.Item(35, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/XyOkWYv31713435061841.jpeg"

# ---- Sheet: 全部类型 - F column updates for unchanged rows 36-49 ----
$ws4.Cells.Item(36, 6).Value = 1328
$ws4.Cells.Item(40, 6).Value = 1418
$ws4.Cells.Item(43, 6).Value = 692
$ws4.Cells.Item(44, 6).Value = 874
$ws4.Cells.Item(48, 6).Value = 261